# Update crypto price/volume data per the Apr 22 2023 GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.294.96'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -4.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.858.31'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -5.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -1.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4500'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -6.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3852'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -5.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.98'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -11.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07884'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -7.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.019'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.35'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.868.42'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -6.88%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.162'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -6.37%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.876'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -5.29%  '
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001031'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.49'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -6.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06534'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.98'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -8.97%  '
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.513'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.296.01'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.75'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -6.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.265'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.086.92'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.80'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.70'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.058'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.461'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -7.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '120.36'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.476'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09283'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9356'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.594'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.279'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02224'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05985'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.21%  '
$ws.Range('E39').Value = '  -3.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.270'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -9.67%  '
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5907'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1879'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.13'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -9.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.253'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5636'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.87'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -9.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.358'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.919'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -7.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06801'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.09'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.06%  '
